# Weekly update: insert two new price rows at the top of the data block
# (rows 1000-1001), pushing the existing rows 1000-1022 down to 1002-1024.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 1000, shifting existing data (old rows
# 1000-1022) down to rows 1002-1024, preserving their formatting/values.
$ws.Rows.Item(1000).Insert()
$ws.Rows.Item(1000).Insert()

# New row 1000: same dimension/category as the (now shifted) row 1002,
# but with this week's updated price figures.
$ws.Cells.Item(1000, 1).Value = 8
$ws.Cells.Item(1000, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1000, 3).Value = "Coquimbo"
$ws.Cells.Item(1000, 4).Value = 45239
$ws.Cells.Item(1000, 5).Value = 4
$ws.Cells.Item(1000, 6).Value = 100112043
$ws.Cells.Item(1000, 7).Value = "Pepino ensalada"
$ws.Cells.Item(1000, 8).Value = "Sin especificar"
$ws.Cells.Item(1000, 9).Value = "Primera"
$ws.Cells.Item(1000, 10).Value = 500
$ws.Cells.Item(1000, 11).Value = 13000
$ws.Cells.Item(1000, 12).Value = 14000
$ws.Cells.Item(1000, 13).Value = 13500
$ws.Cells.Item(1000, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(1000, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1000, 16).Value = 225
$ws.Cells.Item(1000, 17).Value = 60
$ws.Cells.Item(1000, 18).Value = "Hortaliza"

# New row 1001: same week, "Segunda" category.
$ws.Cells.Item(1001, 1).Value = 8
$ws.Cells.Item(1001, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1001, 3).Value = "Coquimbo"
$ws.Cells.Item(1001, 4).Value = 45239
$ws.Cells.Item(1001, 5).Value = 4
$ws.Cells.Item(1001, 6).Value = 100112043
$ws.Cells.Item(1001, 7).Value = "Pepino ensalada"
$ws.Cells.Item(1001, 8).Value = "Sin especificar"
$ws.Cells.Item(1001, 9).Value = "Segunda"
$ws.Cells.Item(1001, 10).Value = 300
$ws.Cells.Item(1001, 11).Value = 9000
$ws.Cells.Item(1001, 12).Value = 10000
$ws.Cells.Item(1001, 13).Value = 9500
$ws.Cells.Item(1001, 14).Value = "$/caja 80 unidades"
$ws.Cells.Item(1001, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1001, 16).Value = 119
$ws.Cells.Item(1001, 17).Value = 80
$ws.Cells.Item(1001, 18).Value = "Hortaliza"
